# Generate Report for Handoff
# Marks the 6daae5bf-... and b1694271-... files as "Ready for handoff" on the
# Overview sheet and on each per-locale sheet (zh-cn, de-de), refreshes the
# "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps, and
# records an Error Detail note explaining that the previous handback was
# against a stale commit.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"

# ---------------------------------------------------------------------------
# Overview sheet: rows 4 (6daae5bf...) and 5 (b1694271...)
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E4").Value = $status
$overview.Range("F4").Value = $status
$overview.Range("G4").Value = "2016-08-20 18:39:01"

$overview.Range("E5").Value = $status
$overview.Range("F5").Value = $status
$overview.Range("G5").Value = "2016-08-20 18:39:01"

# ---------------------------------------------------------------------------
# zh-cn sheet: rows 4 (6daae5bf...) and 5 (b1694271...)
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C4").Value = $status
$zhcn.Range("H4").Value = "2016-08-20 18:38:56"
$zhcn.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a9b4680c280ce0585dd0c37f80bd2e92a934f21/e2e/6daae5bf-6b02-45e6-9a1c-31aacda4d54b.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5d76fad445817b07ec5c052e9f11e0b0e3e839cc/e2e/6daae5bf-6b02-45e6-9a1c-31aacda4d54b.md."

$zhcn.Range("C5").Value = $status
$zhcn.Range("H5").Value = "2016-08-20 18:38:56"
$zhcn.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a9b4680c280ce0585dd0c37f80bd2e92a934f21/e2e/b1694271-d532-4fbc-b30f-dfcb4679f371.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5d76fad445817b07ec5c052e9f11e0b0e3e839cc/e2e/b1694271-d532-4fbc-b30f-dfcb4679f371.md."

# Widen the Error Detail column now that it holds long messages.
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# de-de sheet: rows 4 (6daae5bf...) and 5 (b1694271...)
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C4").Value = $status
$dede.Range("H4").Value = "2016-08-20 18:39:01"
$dede.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a9b4680c280ce0585dd0c37f80bd2e92a934f21/e2e/6daae5bf-6b02-45e6-9a1c-31aacda4d54b.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5d76fad445817b07ec5c052e9f11e0b0e3e839cc/e2e/6daae5bf-6b02-45e6-9a1c-31aacda4d54b.md."

$dede.Range("C5").Value = $status
$dede.Range("H5").Value = "2016-08-20 18:39:01"
$dede.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2a9b4680c280ce0585dd0c37f80bd2e92a934f21/e2e/b1694271-d532-4fbc-b30f-dfcb4679f371.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5d76fad445817b07ec5c052e9f11e0b0e3e839cc/e2e/b1694271-d532-4fbc-b30f-dfcb4679f371.md."

# Widen the Error Detail column now that it holds long messages.
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
